$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1129.8064
$ws.Range("I15").Value = 1129.8064
$ws.Range("K15").Value = 3389.4192
$ws.Range("M15").Value = -3220.4192
$ws.Range("H96").Value = 1715.762
$ws.Range("I96").Value = 449.22223
$ws.Range("J96").Value = 2665.6667
$ws.Range("K96").Value = 1347.66669
$ws.Range("L96").Value = 7997.000100000001
$ws.Range("M96").Value = 25.33330999999998
$ws.Range("N96").Value = -10743.0001
$ws.Range("H106").Value = 5330.304
$ws.Range("I106").Value = 5390.773
$ws.Range("K106").Value = 5390.773
$ws.Range("M106").Value = -4759.773
$ws.Range("H132").Value = 199999.5
$ws.Range("I132").Value = 199999.5
$ws.Range("K132").Value = 599998.5
$ws.Range("M132").Value = -597468.5
$ws.Range("H138").Value = 4432.3823
$ws.Range("J138").Value = 4833.3335
$ws.Range("L138").Value = 14500.0005
$ws.Range("N138").Value = -24780.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 598
$ws.Range("I25").Value = 547
$ws.Range("J25").Value = 700
$ws.Range("K25").Value = 547
$ws.Range("L25").Value = 700
$ws.Range("M25").Value = -145
$ws.Range("N25").Value = -1504
$ws.Range("H61").Value = 2113
$ws.Range("I61").Value = 2111.4614
$ws.Range("J61").Value = 2133
$ws.Range("K61").Value = 2111.4614
$ws.Range("L61").Value = 2133
$ws.Range("M61").Value = -1899.4614
$ws.Range("N61").Value = -2557
$ws.Range("H97").Value = 813.5454999999999
$ws.Range("J97").Value = 3200
$ws.Range("L97").Value = 3200
$ws.Range("N97").Value = -4192
$ws.Range("H132").Value = 5775.6553
$ws.Range("I132").Value = 3718.7568
$ws.Range("K132").Value = 11156.2704
$ws.Range("M132").Value = -8626.270400000001
$ws.Range("H136").Value = 2113
$ws.Range("I136").Value = 2111.4614
$ws.Range("J136").Value = 2133
$ws.Range("K136").Value = 6334.3842
$ws.Range("L136").Value = 6399
$ws.Range("M136").Value = -3784.3842
$ws.Range("N136").Value = -11499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 53300.65
$ws.Range("I20").Value = 74916.28999999999
$ws.Range("K20").Value = 74916.28999999999
$ws.Range("M20").Value = -74669.28999999999
$ws.Range("H80").Value = 559.4167
$ws.Range("I80").Value = 782.5
$ws.Range("J80").Value = 400.07144
$ws.Range("K80").Value = 782.5
$ws.Range("L80").Value = 400.07144
$ws.Range("M80").Value = 215.5
$ws.Range("N80").Value = -2396.07144
$ws.Range("H83").Value = 559.4167
$ws.Range("I83").Value = 782.5
$ws.Range("J83").Value = 400.07144
$ws.Range("K83").Value = 3912.5
$ws.Range("L83").Value = 2000.3572
$ws.Range("M83").Value = 1079.5
$ws.Range("N83").Value = -11984.3572
$ws.Range("H99").Value = 11279.6
$ws.Range("I99").Value = 12366.223
$ws.Range("K99").Value = 12366.223
$ws.Range("M99").Value = -10868.223
$ws.Range("H105").Value = 4077.4167
$ws.Range("I105").Value = 1774.7894
$ws.Range("K105").Value = 1774.7894
$ws.Range("M105").Value = -27.78939999999989
$ws.Range("H134").Value = 2823.1538
$ws.Range("I134").Value = 2445.1
$ws.Range("K134").Value = 7335.299999999999
$ws.Range("M134").Value = -4800.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 4500
$ws.Range("J38").Value = 4500
$ws.Range("L38").Value = 4500
$ws.Range("N38").Value = -5254
$ws.Range("H46").Value = 4500
$ws.Range("J46").Value = 4500
$ws.Range("L46").Value = 4500
$ws.Range("N46").Value = -4922
$ws.Range("H55").Value = 8944
$ws.Range("I55").Value = 10258.667
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 10258.667
$ws.Range("L55").Value = 5000
$ws.Range("M55").Value = -9943.666999999999
$ws.Range("N55").Value = -5630
$ws.Range("H94").Value = 5568.095
$ws.Range("I94").Value = 14604.714
$ws.Range("K94").Value = 14604.714
$ws.Range("M94").Value = -14153.714
$ws.Range("H107").Value = 1811.6875
$ws.Range("I107").Value = 1549.25
$ws.Range("K107").Value = 1549.25
$ws.Range("M107").Value = 370.75
$ws.Range("H134").Value = 2459.4827
$ws.Range("I134").Value = 2060.682
$ws.Range("J134").Value = 3712.8572
$ws.Range("K134").Value = 6182.045999999999
$ws.Range("L134").Value = 11138.5716
$ws.Range("M134").Value = -3647.045999999999
$ws.Range("N134").Value = -16208.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1501.7273
$ws.Range("I33").Value = 110.6
$ws.Range("J33").Value = 2661
$ws.Range("K33").Value = 663.5999999999999
$ws.Range("L33").Value = 15966
$ws.Range("M33").Value = -380.5999999999999
$ws.Range("N33").Value = -16532
$ws.Range("H68").Value = 5085.5
$ws.Range("I68").Value = 2324.5
$ws.Range("K68").Value = 6973.5
$ws.Range("M68").Value = -6162.5
$ws.Range("H71").Value = 5085.5
$ws.Range("I71").Value = 2324.5
$ws.Range("K71").Value = 20920.5
$ws.Range("M71").Value = -16864.5
$ws.Range("H98").Value = 510.42856
$ws.Range("I98").Value = 548.8
$ws.Range("J98").Value = 489.1111
$ws.Range("K98").Value = 1646.4
$ws.Range("L98").Value = 1467.3333
$ws.Range("M98").Value = -148.3999999999999
$ws.Range("N98").Value = -4463.3333
$ws.Range("H113").Value = 29309
$ws.Range("J113").Value = 36082.47
$ws.Range("L113").Value = 108247.41
$ws.Range("N113").Value = -112587.41
$ws.Range("H137").Value = 10503.223
$ws.Range("J137").Value = 11187.375
$ws.Range("L137").Value = 33562.125
$ws.Range("N137").Value = -43762.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 11166.333
$ws.Range("I18").Value = 16000
$ws.Range("J18").Value = 8749.5
$ws.Range("K18").Value = 16000
$ws.Range("L18").Value = 8749.5
$ws.Range("M18").Value = -15707
$ws.Range("N18").Value = -9335.5
$ws.Range("H113").Value = 1853.3103
$ws.Range("I113").Value = 1815.3914
$ws.Range("J113").Value = 1998.6666
$ws.Range("K113").Value = 1815.3914
$ws.Range("L113").Value = 1998.6666
$ws.Range("M113").Value = 354.6086
$ws.Range("N113").Value = -6338.6666
$ws.Range("H132").Value = 24312.438
$ws.Range("I132").Value = 36429.9
$ws.Range("J132").Value = 4116.6665
$ws.Range("K132").Value = 109289.7
$ws.Range("L132").Value = 12349.9995
$ws.Range("M132").Value = -106759.7
$ws.Range("N132").Value = -17409.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4371.48
$ws.Range("I46").Value = 1283.8
$ws.Range("J46").Value = 5143.4
$ws.Range("K46").Value = 1283.8
$ws.Range("L46").Value = 5143.4
$ws.Range("M46").Value = -1095.8
$ws.Range("N46").Value = -5519.4
$ws.Range("H55").Value = 1437.9688
$ws.Range("I55").Value = 1272.3846
$ws.Range("K55").Value = 1272.3846
$ws.Range("M55").Value = -1099.3846
$ws.Range("H61").Value = 9887.058999999999
$ws.Range("I61").Value = 9720.357
$ws.Range("J61").Value = 10665
$ws.Range("K61").Value = 9720.357
$ws.Range("L61").Value = 10665
$ws.Range("M61").Value = -9518.357
$ws.Range("N61").Value = -11069
$ws.Range("H106").Value = 24833
$ws.Range("J106").Value = 24833
$ws.Range("L106").Value = 24833
$ws.Range("N106").Value = -27357
$ws.Range("H113").Value = 9887.058999999999
$ws.Range("I113").Value = 9720.357
$ws.Range("J113").Value = 10665
$ws.Range("K113").Value = 9720.357
$ws.Range("L113").Value = 10665
$ws.Range("M113").Value = -7550.357
$ws.Range("N113").Value = -15005
$ws.Range("H122").Value = 4222.5557
$ws.Range("I122").Value = 3940.8
$ws.Range("J122").Value = 4574.75
$ws.Range("K122").Value = 11822.4
$ws.Range("L122").Value = 13724.25
$ws.Range("M122").Value = -9372.400000000001
$ws.Range("N122").Value = -18624.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 10000
$ws.Range("I101").Value = 10000
$ws.Range("K101").Value = 10000
$ws.Range("M101").Value = -6755
$ws.Range("H113").Value = 11663.333
$ws.Range("I113").Value = 9995
$ws.Range("J113").Value = 15000
$ws.Range("K113").Value = 29985
$ws.Range("L113").Value = 45000
$ws.Range("M113").Value = -27815
$ws.Range("N113").Value = -49340
$ws.Range("H122").Value = 41227.1
$ws.Range("I122").Value = 2443.762
$ws.Range("K122").Value = 7331.286
$ws.Range("M122").Value = -4881.286
$ws.Range("H123").Value = 59999
$ws.Range("J123").Value = 59999
$ws.Range("L123").Value = 59999
$ws.Range("N123").Value = -69799
